# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.064.18"
$ws.Range("E2").Value = "  -3.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.13"
$ws.Range("E3").Value = "  -3.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -1.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.06"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("E6").Value = "  -1.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4680"
$ws.Range("E7").Value = "  -5.96%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4011"
$ws.Range("E8").Value = "  -4.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.11"
$ws.Range("E9").Value = "  -2.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08350"
$ws.Range("E10").Value = "  -10.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.041"
$ws.Range("E11").Value = "  -4.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.09"
$ws.Range("E12").Value = "  -2.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.921.97"
$ws.Range("E13").Value = "  -2.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.438"
$ws.Range("E14").Value = "  -5.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.049"
$ws.Range("E15").Value = "  -6.21%  "

$ws.Range("E16").Value = "  -1.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.32"
$ws.Range("E17").Value = "  -2.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001061"
$ws.Range("E18").Value = "  -4.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06573"
$ws.Range("E19").Value = "  -2.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.93"
$ws.Range("E20").Value = "  -6.29%  "

$ws.Range("E21").Value = "  -1.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.703"
$ws.Range("E22").Value = "  -4.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.059.19"
$ws.Range("E23").Value = "  -3.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("E24").Value = "  -5.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.280"
$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.115.00"
$ws.Range("E26").Value = "  -4.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.92"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.97"
$ws.Range("E28").Value = "  -3.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.125"
$ws.Range("E29").Value = "  -5.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.652"
$ws.Range("E30").Value = "  -9.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.78"
$ws.Range("E31").Value = "  -3.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9704"
$ws.Range("E32").Value = "  -6.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09547"
$ws.Range("E33").Value = "  -2.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.440"
$ws.Range("E34").Value = "  -3.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.632"
$ws.Range("E35").Value = "  -3.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.521"
$ws.Range("E36").Value = "  -4.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.758"
$ws.Range("E37").Value = "  -2.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02295"
$ws.Range("E38").Value = "  -5.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06124"
$ws.Range("E39").Value = "  -4.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.213"
$ws.Range("E40").Value = "  -8.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6114"
$ws.Range("E41").Value = "  -5.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.99"
$ws.Range("E42").Value = "  -3.90%  "

$ws.Range("E43").Value = "  -1.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1896"
$ws.Range("E44").Value = "  -5.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.301"
$ws.Range("E45").Value = "  -3.57%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.77"
$ws.Range("E46").Value = "  -3.73%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5824"
$ws.Range("E47").Value = "  -5.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.015"
$ws.Range("E48").Value = "  -7.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.446"
$ws.Range("E49").Value = "  -1.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06839"
$ws.Range("E50").Value = "  -1.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.06"
$ws.Range("E51").Value = "  -3.62%  "
